$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Score summary block (rows 10-12) -------------------------------------
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("B10").Value = 12
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 13
$ws.Range("E10").Value = 28

$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("B12").Value = 48
$ws.Range("C12").Value = -3
$ws.Range("E12").Value = "45/112"

# --- Drop the 2nd/3rd "Student Ans / Correct Ans" repeated blocks ---------
# Third block (columns G:H) is removed entirely for rows 15-40.
$ws.Range("G15:H40").Delete()
# Second block (columns D:E) only keeps its header + first two answer rows.
$ws.Range("D19:E40").Delete()

# --- Fill in the student answers for the first block (column A) ----------
$ws.Range("D16").Style = "correctStyle"
$ws.Range("D16").Value = "Option A"

$ws.Range("A17").Style = "correctStyle"
$ws.Range("A17").Value = "Option D"

$ws.Range("A18").Style = "correctStyle"
$ws.Range("A18").Value = "Option B"

$ws.Range("A20").Style = "correctStyle"
$ws.Range("A20").Value = "Option B"

$ws.Range("A21").Style = "correctStyle"
$ws.Range("A21").Value = "Option C"

$ws.Range("A26").Style = "incorrectStyle"
$ws.Range("A26").Value = "Option B"

$ws.Range("A27").Style = "correctStyle"
$ws.Range("A27").Value = "Option A"

$ws.Range("A29").Style = "correctStyle"
$ws.Range("A29").Value = "Option D"

$ws.Range("A30").Style = "correctStyle"
$ws.Range("A30").Value = "Option B"

$ws.Range("A32").Style = "correctStyle"
$ws.Range("A32").Value = "Option C"

$ws.Range("A34").Style = "incorrectStyle"
$ws.Range("A34").Value = "Option A"

$ws.Range("A35").Style = "correctStyle"
$ws.Range("A35").Value = "Option D"

$ws.Range("A36").Style = "correctStyle"
$ws.Range("A36").Value = "Option A"

$ws.Range("A37").Style = "incorrectStyle"
$ws.Range("A37").Value = "Option B"

$ws.Range("A39").Style = "correctStyle"
$ws.Range("A39").Value = "Option D"
